{"js": "// Translate the English subtitle text to Swahili.\n// Each pair is [English source text, Swahili replacement text].\nconst replacements = [\n  [\"Format has been corrected not the timing\",\n   \"Umbizo limesahihishwa sio wakati\"],\n  [\"I added 25 seconds to each timing to correct for the intro song -john argentino\",\n   \"Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino\"],\n  [\"The airport problem - subtitles:\",\n   \"Tatizo la uwanja wa ndege - manukuu:\"],\n  [\"The administrations of three\",\n   \"Utawala wa tatu\"],\n  [\"neighboring cities: A, B and C decided\",\n   \"miji jirani: A, B na C waliamua\"],\n  [\"to build an airport dividing the costs of\",\n   \"kujenga uwanja wa ndege unaogawanya gharama za\"],\n  [\"implementation. The condition on the\",\n   \"utekelezaji. Hali juu ya\"],\n  [\"choice of the most suitable place is\",\n   \"uchaguzi wa mahali pa kufaa zaidi ni\"],\n  [\"that the sum of the distances from each\",\n   \"kwamba jumla ya umbali kutoka kwa kila mmoja\"],\n  [\"city to the airport is as small as\",\n   \"mji kwa uwanja wa ndege ni ndogo kama\"],\n  [\"possible. The team of experts in charge\",\n   \"inawezekana. Timu ya wataalam wanaohusika\"],\n  [\"of the work has created a model to get\",\n   \"ya kazi imeunda mfano wa kupata\"],\n  [\"a preliminary idea of where to place the\",\n   \"wazo la awali la mahali pa kuweka\"],\n  [\"structure. At their disposal there are\",\n   \"muundo. At their disposal there are\"],\n  [\"[Music]\",\n   \"[Muziki]\"],\n];\n\nconst body = context.document.body;\n\nfor (const [source, target] of replacements) {\n  const results = body.search(source, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(target, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English subtitle text to Swahili using Find & Replace.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Format has been corrected not the timing\", \"Umbizo limesahihishwa sio wakati\"),\n    @(\"I added 25 seconds to each timing to correct for the intro song -john argentino\", \"Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino\"),\n    @(\"The airport problem - subtitles:\", \"Tatizo la uwanja wa ndege - manukuu:\"),\n    @(\"The administrations of three\", \"Utawala wa tatu\"),\n    @(\"neighboring cities: A, B and C decided\", \"miji jirani: A, B na C waliamua\"),\n    @(\"to build an airport dividing the costs of\", \"kujenga uwanja wa ndege unaogawanya gharama za\"),\n    @(\"implementation. The condition on the\", \"utekelezaji. Hali juu ya\"),\n    @(\"choice of the most suitable place is\", \"uchaguzi wa mahali pa kufaa zaidi ni\"),\n    @(\"that the sum of the distances from each\", \"kwamba jumla ya umbali kutoka kwa kila mmoja\"),\n    @(\"city to the airport is as small as\", \"mji kwa uwanja wa ndege ni ndogo kama\"),\n    @(\"possible. The team of experts in charge\", \"inawezekana. Timu ya wataalam wanaohusika\"),\n    @(\"of the work has created a model to get\", \"ya kazi imeunda mfano wa kupata\"),\n    @(\"a preliminary idea of where to place the\", \"wazo la awali la mahali pa kuweka\"),\n    @(\"structure. At their disposal there are\", \"muundo. At their disposal there are\"),\n    @(\"[Music]\", \"[Muziki]\")\n)\n\nforeach ($pair in $replacements) {\n    $source = $pair[0]\n    $target = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $source\n    $find.Replacement.Text = $target\n    $find.Execute($source, $false, $false, $false, $false, $false, $true, 1, $false, $target, 2)\n}\n"}
